# Removed hard coded data for the products tests:
# Populate the ProductsTests sheet with the new TC_009 / TC_010 rows
# (expected cart contents, sort order) instead of relying on the
# hard-coded single-row fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductsTests")

# New cell values - order chosen to match the original shared-string
# insertion order (TC ids first, then the long product/sort strings).
$ws.Range("A2").Value = "TC_009"
$ws.Range("A3").Value = "TC_010"
$ws.Range("C1").Value = "Sauce Labs Backpack, Sauce Labs Bike Light, Sauce Labs Bolt T-Shirt, Sauce Labs Fleece Jacket, Sauce Labs Onesie, Test.allTheThings() T-Shirt (Red)"
$ws.Range("B2").Value = "Test.allTheThings() T-Shirt (Red), Sauce Labs Onesie, Sauce Labs Fleece Jacket, Sauce Labs Bolt T-Shirt, Sauce Labs Bike Light, Sauce Labs Backpack"
$ws.Range("B3").Value = "Name (A to Z), Name (Z to A), Price (low to high), Price (high to low)"

# Wrap the long text cells so they are readable
$ws.Range("C1").WrapText = $true
$ws.Range("B2").WrapText = $true

# Row heights: row 1 grows to fit two wrapped lines, row 2 is sized
# generously for its longer wrapped text.
$ws.Rows.Item(1).RowHeight = 28.8
$ws.Rows.Item(2).RowHeight = 49.8

# Column widths for the newly-populated / widened columns
$ws.Columns.Item(2).ColumnWidth = 56.21875
$ws.Columns.Item(3).ColumnWidth = 78.88671875

# Leave the selection where the user ended up after entering the data
$ws.Range("B11").Select()
